$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the data area so that cells which become
# empty in the new layout (e.g. old column G values on "regular" rows)
# do not linger around.
$ws.UsedRange.ClearContents() | Out-Null

# Title
$ws.Range("A1").Value = "Model training tracking"

# Header row (row 3) - two new columns were inserted:
#   "zeroshot minsim" (new col G) and "seed words included" (new col J)
$ws.Range("B3").Value = "split_comp"
$ws.Range("C3").Value = "split"
$ws.Range("D3").Value = "sample"
$ws.Range("E3").Value = "min cluster"
$ws.Range("F3").Value = "n components"
$ws.Range("G3").Value = "zeroshot minsim"
$ws.Range("H3").Value = "training time"
$ws.Range("I3").Value = "nr clusters"
$ws.Range("J3").Value = "seed words included"
$ws.Range("K3").Value = "topics produced"
$ws.Range("L3").Value = "topic quality (eigene Beurteilung)"

# Row 4 - regular / full random / sentences
$ws.Range("A4").Value = "regular"
$ws.Range("B4").Value = "full random"
$ws.Range("C4").Value = "sentences"
$ws.Range("D4").Value = 3000
$ws.Range("E4").Value = 200
$ws.Range("F4").Value = 5
$ws.Range("H4").Value = "2000s (ca)"
$ws.Range("I4").Value = "auto"
$ws.Range("J4").Value = "no"
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = "medium"

# Row 5 - regular / full random / sentences
$ws.Range("A5").Value = "regular"
$ws.Range("B5").Value = "full random"
$ws.Range("C5").Value = "sentences"
$ws.Range("D5").Value = 3000
$ws.Range("E5").Value = 200
$ws.Range("F5").Value = 10
$ws.Range("H5").Value = "2000s (ca)"
$ws.Range("I5").Value = "auto"
$ws.Range("J5").Value = "no"
$ws.Range("K5").Value = 8
$ws.Range("L5").Value = "bad"

# Row 6 - zeroshot / full random / sentences
$ws.Range("A6").Value = "zeroshot"
$ws.Range("B6").Value = "full random"
$ws.Range("C6").Value = "sentences"
$ws.Range("D6").Value = 5000
$ws.Range("E6").Value = 200
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 0.7
$ws.Range("H6").Value = "15000s"
$ws.Range("I6").Value = "auto!?"
$ws.Range("J6").Value = "no"
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = "bad"

# Row 7 (new) - zeroshot / full random / sentences
$ws.Range("A7").Value = "zeroshot"
$ws.Range("B7").Value = "full random"
$ws.Range("C7").Value = "sentences"
$ws.Range("D7").Value = 3000
$ws.Range("E7").Value = 200
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 0.3
$ws.Range("H7").Value = "3000s"
$ws.Range("I7").Value = "null"
$ws.Range("J7").Value = "no"
$ws.Range("K7").Value = "100+"
$ws.Range("L7").Value = "quite good"

# Row 8 (new) - regular / full random / sentences
$ws.Range("A8").Value = "regular"
$ws.Range("B8").Value = "full random"
$ws.Range("C8").Value = "sentences"
$ws.Range("D8").Value = 3000
$ws.Range("E8").Value = 200
$ws.Range("F8").Value = 5
$ws.Range("H8").Value = "2000s (ca)"
$ws.Range("I8").Value = "auto"
$ws.Range("J8").Value = "yes"
$ws.Range("K8").Value = 50
$ws.Range("L8").Value = "good"

# Re-fit the two newly introduced columns to their content, matching the
# author's column-width adjustments for the inserted columns.
$ws.Columns("J:K").AutoFit() | Out-Null

# Restore the last active selection as recorded in the saved file.
$ws.Range("J12").Select() | Out-Null
